# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets,
# matching the output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1) - row => new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 167
$ws1.Range("F3").Value = 157
$ws1.Range("F4").Value = 162
$ws1.Range("F5").Value = 4698
$ws1.Range("F8").Value = 519
$ws1.Range("F9").Value = 477
$ws1.Range("F10").Value = 23
$ws1.Range("F12").Value = 1339
$ws1.Range("F13").Value = 2805
$ws1.Range("F14").Value = 385
$ws1.Range("F15").Value = 92
$ws1.Range("F16").Value = 79
$ws1.Range("F17").Value = 67
$ws1.Range("F18").Value = 2365
$ws1.Range("F19").Value = 104
$ws1.Range("F20").Value = 77
$ws1.Range("F22").Value = 165
$ws1.Range("F23").Value = 113
$ws1.Range("F25").Value = 233
$ws1.Range("F26").Value = 41

# Sheet "全部类型" (sheetId 4) - row => new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 167
$ws4.Range("F3").Value = 157
$ws4.Range("F4").Value = 162
$ws4.Range("F6").Value = 4698
$ws4.Range("F9").Value = 519
$ws4.Range("F10").Value = 477
$ws4.Range("F11").Value = 23
$ws4.Range("F13").Value = 1339
$ws4.Range("F14").Value = 2805
$ws4.Range("F15").Value = 385
$ws4.Range("F16").Value = 92
$ws4.Range("F17").Value = 79
$ws4.Range("F18").Value = 67
$ws4.Range("F19").Value = 2365
$ws4.Range("F20").Value = 104
$ws4.Range("F21").Value = 77
$ws4.Range("F23").Value = 165
$ws4.Range("F24").Value = 113
$ws4.Range("F26").Value = 233
$ws4.Range("F27").Value = 41
